$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look like pure numbers,
# so Excel keeps them as text (matching the source data which stores prices as strings).
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D14", "D15", "D18", "D19", "D20", "D21", "D23", "D24", "D26", "D27", "D29", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D47", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '59.476.55'
$ws.Range('E2').Value = '  +2.86%  '
$ws.Range('D3').Value = '3.188.22'
$ws.Range('E3').Value = '  +1.91%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '533.61'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').Value = '142.56'
$ws.Range('E6').Value = '  +2.65%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('E8').Value = '  +10.82%  '
$ws.Range('D9').Value = '7.30'
$ws.Range('E9').Value = '  -0.30%  '
$ws.Range('D10').Value = '0.439'
$ws.Range('E10').Value = '  +6.64%  '
$ws.Range('D11').Value = '0.112'
$ws.Range('E11').Value = '  +4.36%  '
$ws.Range('D12').Value = '3.740.14'
$ws.Range('E12').Value = '  +1.98%  '
$ws.Range('E13').Value = '  +1.70%  '
$ws.Range('D14').Value = '25.90'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').Value = '0.0000172'
$ws.Range('E15').Value = '  +4.97%  '
$ws.Range('D16').Value = '59.622.18'
$ws.Range('E16').Value = '  +2.93%  '
$ws.Range('D17').Value = '3.200.09'
$ws.Range('E17').Value = '  +2.30%  '
$ws.Range('D18').Value = '6.27'
$ws.Range('E18').Value = '  +2.90%  '
$ws.Range('D19').Value = '13.04'
$ws.Range('E19').Value = '  +2.71%  '
$ws.Range('D20').Value = '8.21'
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('D21').Value = '376.36'
$ws.Range('E21').Value = '  +2.36%  '
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').Value = '0.532'
$ws.Range('E23').Value = '  +5.12%  '
$ws.Range('D24').Value = '70.21'
$ws.Range('E24').Value = '  +1.37%  '
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('D27').Value = '8.41'
$ws.Range('E27').Value = '  +14.89%  '
$ws.Range('D28').Value = '0.0₃0879'
$ws.Range('E28').Value = '  +1.74%  '
$ws.Range('D29').Value = '22.45'
$ws.Range('E29').Value = '  +4.66%  '
$ws.Range('E30').Value = '  +1.18%  '
$ws.Range('D31').Value = '6.12'
$ws.Range('E31').Value = '  +0.22%  '
$ws.Range('D32').Value = '5.29'
$ws.Range('E32').Value = '  +2.45%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').Value = '6.38'
$ws.Range('E34').Value = '  +4.76%  '
$ws.Range('D35').Value = '157.27'
$ws.Range('E35').Value = '  -1.36%  '
$ws.Range('D36').Value = '1.34'
$ws.Range('E36').Value = '  +3.68%  '
$ws.Range('B37').Value = 'EnergySwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D37').Value = '25.65'
$ws.Range('E37').Value = '  +0.82%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.0715'
$ws.Range('E38').Value = '  +6.45%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '2.716.55'
$ws.Range('E39').Value = '  +7.43%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '1.72'
$ws.Range('E40').Value = '  +2.55%  '
$ws.Range('D41').Value = '4.29'
$ws.Range('E41').Value = '  +4.66%  '
$ws.Range('D42').Value = '0.728'
$ws.Range('E42').Value = '  +4.14%  '
$ws.Range('D43').Value = '0.0293'
$ws.Range('E43').Value = '  +8.50%  '
$ws.Range('D44').Value = '39.27'
$ws.Range('E44').Value = '  +3.88%  '
$ws.Range('E45').Value = '  +0.15%  '
$ws.Range('D46').Value = '3.231.97'
$ws.Range('E46').Value = '  +1.93%  '
$ws.Range('D47').Value = '0.994'
$ws.Range('E47').Value = '  +1.45%  '
$ws.Range('E48').Value = '  +11.77%  '
$ws.Range('D49').Value = '6.23'
$ws.Range('E49').Value = '  +1.54%  '
$ws.Range('D50').Value = '20.54'
$ws.Range('E50').Value = '  +3.93%  '
$ws.Range('D51').Value = '0.764'
$ws.Range('E51').Value = '  +2.37%  '
